# Re-applies the table style used by the table on the "B1 - Types of
# financial documents" slide: it was switched from the default
# "Table_0" style ({1DF7846D-7C2D-4AA4-89BC-015AC81FCC01}) to the
# built-in style {45C1CBD5-C263-4CD8-9565-12E2DE40B586} via the Table
# Design ribbon ("Table Styles" gallery).

$p = $ppt.ActivePresentation

$targetStyleId = "{45C1CBD5-C263-4CD8-9565-12E2DE40B586}"

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $s = $p.Slides.Item($i)
    for ($j = 1; $j -le $s.Shapes.Count; $j++) {
        $sh = $s.Shapes.Item($j)
        if ($sh.HasTable) {
            $tbl = $sh.Table
            $tbl.ApplyStyle($targetStyleId)
        }
    }
}
